# Applies the cryptos.xlsx crypto-price refresh described in the commit
# 'Updated cryptos list on Thu Sep 19 17:37:38 UTC 2024 with GitHub Actions'.
# Every data cell (rows 2-51, cols D/E, plus the B/C/D/E swap of rows 30<->31
# and 50<->51) is rewritten as literal text, matching the original inlineStr
# cell contents (prices use '.' as a thousands separator, so they must stay
# text rather than become numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Writing a number-looking string via .Value would make Excel auto-convert it
# to a numeric value (e.g. '1.00' -> 1). Prefixing with an apostrophe forces
# text entry (like typing '0.999 into a cell); resetting the cell .Style to
# 'Normal' afterwards drops the transient quotePrefix formatting that the
# apostrophe trick applies, so the cell style stays identical to the original.
function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.Value = "'" + $val
    $r.Style = 'Normal'
}

Set-TextCell 'D2' '63.476.46'
Set-TextCell 'E2' '  +5.94%  '
Set-TextCell 'D3' '2.462.15'
Set-TextCell 'E3' '  +6.86%  '
Set-TextCell 'D4' '0.999'
Set-TextCell 'E4' '  +0.01%  '
Set-TextCell 'D5' '569.13'
Set-TextCell 'E5' '  +5.11%  '
Set-TextCell 'D6' '142.86'
Set-TextCell 'E6' '  +10.78%  '
Set-TextCell 'E7' '  +0.04%  '
Set-TextCell 'D8' '0.591'
Set-TextCell 'E8' '  +3.87%  '
Set-TextCell 'D9' '2.459.43'
Set-TextCell 'E9' '  +6.90%  '
Set-TextCell 'E10' '  +5.00%  '
Set-TextCell 'D11' '5.75'
Set-TextCell 'E11' '  +4.75%  '
Set-TextCell 'E13' '  +7.22%  '
Set-TextCell 'D14' '26.41'
Set-TextCell 'E14' '  +14.11%  '
Set-TextCell 'D15' '2.896.40'
Set-TextCell 'E15' '  +6.94%  '
Set-TextCell 'D16' '63.347.91'
Set-TextCell 'E16' '  +5.91%  '
Set-TextCell 'E17' '  +8.95%  '
Set-TextCell 'D18' '2.460.35'
Set-TextCell 'E18' '  +6.84%  '
Set-TextCell 'D19' '11.30'
Set-TextCell 'E19' '  +8.54%  '
Set-TextCell 'D20' '342.96'
Set-TextCell 'E20' '  +10.29%  '
Set-TextCell 'E21' '  +7.63%  '
Set-TextCell 'D22' '6.83'
Set-TextCell 'E22' '  +5.38%  '
Set-TextCell 'E23' '  -0.15%  '
Set-TextCell 'D24' '65.76'
Set-TextCell 'E24' '  +3.35%  '
Set-TextCell 'E25' '  +3.93%  '
Set-TextCell 'D26' '1.00'
Set-TextCell 'E26' '  +0.12%  '
Set-TextCell 'D27' '1.54'
Set-TextCell 'E27' '  +15.10%  '
Set-TextCell 'D28' '8.23'
Set-TextCell 'E28' '  +6.95%  '
Set-TextCell 'E29' '  +13.13%  '
Set-TextCell 'B30' 'PEPE'
Set-TextCell 'C30' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextCell 'D30' '0.0₃0822'
Set-TextCell 'E30' '  +14.94%  '
Set-TextCell 'B31' 'Aptos'
Set-TextCell 'C31' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextCell 'D31' '6.89'
Set-TextCell 'E31' '  +18.84%  '
Set-TextCell 'D32' '1.86'
Set-TextCell 'E32' '  +9.48%  '
Set-TextCell 'D33' '175.19'
Set-TextCell 'E33' '  +2.32%  '
Set-TextCell 'E34' '  +11.58%  '
Set-TextCell 'E35' '  +6.11%  '
Set-TextCell 'E36' '  +6.84%  '
Set-TextCell 'D37' '371.25'
Set-TextCell 'E37' '  +19.15%  '
Set-TextCell 'E38' '  +10.85%  '
Set-TextCell 'E39' '  +0.02%  '
Set-TextCell 'D40' '1.74'
Set-TextCell 'E40' '  +15.82%  '
Set-TextCell 'D41' '0.999'
Set-TextCell 'E41' '  +0.09%  '
Set-TextCell 'D42' '40.39'
Set-TextCell 'E42' '  +6.45%  '
Set-TextCell 'D43' '150.33'
Set-TextCell 'E43' '  +10.23%  '
Set-TextCell 'E44' '  +9.71%  '
Set-TextCell 'D45' '20.75'
Set-TextCell 'E45' '  +12.33%  '
Set-TextCell 'D46' '0.599'
Set-TextCell 'E46' '  +6.08%  '
Set-TextCell 'D47' '0.0967'
Set-TextCell 'E47' '  +3.37%  '
Set-TextCell 'D48' '0.0526'
Set-TextCell 'E48' '  +8.06%  '
Set-TextCell 'E49' '  +7.03%  '
Set-TextCell 'B50' 'EnergySwap'
Set-TextCell 'C50' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D50' '18.19'
Set-TextCell 'E50' '  +9.56%  '
Set-TextCell 'B51' 'BabyDogeCoin'
Set-TextCell 'C51' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell 'D51' '0.0₆0230'
Set-TextCell 'E51' '  +5.35%  '
